$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(62, 8).Value = 6950099.5
$ws.Cells.Item(62, 9).Value = 15628725
$ws.Cells.Item(62, 10).Value = 7198.6
$ws.Cells.Item(62, 11).Value = 15628725
$ws.Cells.Item(62, 12).Value = 7198.6
$ws.Cells.Item(62, 13).Value = -15628101
$ws.Cells.Item(62, 14).Value = -8446.6
$ws.Cells.Item(65, 8).Value = 6950099.5
$ws.Cells.Item(65, 9).Value = 15628725
$ws.Cells.Item(65, 10).Value = 7198.6
$ws.Cells.Item(65, 11).Value = 78143625
$ws.Cells.Item(65, 12).Value = 35993
$ws.Cells.Item(65, 13).Value = -78140505
$ws.Cells.Item(65, 14).Value = -42233
$ws.Cells.Item(92, 8).Value = 653.4
$ws.Cells.Item(92, 9).Value = 656.125
$ws.Cells.Item(92, 10).Value = 642.5
$ws.Cells.Item(92, 11).Value = 656.125
$ws.Cells.Item(92, 12).Value = 642.5
$ws.Cells.Item(92, 13).Value = 591.875
$ws.Cells.Item(92, 14).Value = -3138.5
$ws.Cells.Item(97, 8).Value = 4832.6665
$ws.Cells.Item(97, 10).Value = 4832.6665
$ws.Cells.Item(97, 12).Value = 14497.9995
$ws.Cells.Item(97, 14).Value = -15489.9995
$ws.Cells.Item(98, 8).Value = 1863.6086
$ws.Cells.Item(98, 9).Value = 1450.7894
$ws.Cells.Item(98, 10).Value = 3824.5
$ws.Cells.Item(98, 11).Value = 1450.7894
$ws.Cells.Item(98, 12).Value = 3824.5
$ws.Cells.Item(98, 13).Value = 47.21060000000011
$ws.Cells.Item(98, 14).Value = -6820.5
$ws.Cells.Item(107, 8).Value = 48645.855
$ws.Cells.Item(107, 9).Value = 56653.832
$ws.Cells.Item(107, 11).Value = 56653.832
$ws.Cells.Item(107, 13).Value = -54733.832
$ws.Cells.Item(110, 8).Value = 42189
$ws.Cells.Item(110, 10).Value = 42189
$ws.Cells.Item(110, 12).Value = 42189
$ws.Cells.Item(110, 14).Value = -50369
$ws.Cells.Item(111, 8).Value = 79143.08
$ws.Cells.Item(111, 10).Value = 2916
$ws.Cells.Item(111, 12).Value = 8748
$ws.Cells.Item(111, 14).Value = -14882
$ws.Cells.Item(112, 8).Value = 1983.3334
$ws.Cells.Item(112, 10).Value = 1983.3334
$ws.Cells.Item(112, 12).Value = 5950.0002
$ws.Cells.Item(112, 14).Value = -8166.0002
$ws.Cells.Item(116, 8).Value = 30956.625
$ws.Cells.Item(116, 9).Value = 30912
$ws.Cells.Item(116, 10).Value = 31001.25
$ws.Cells.Item(116, 11).Value = 30912
$ws.Cells.Item(116, 12).Value = 31001.25
$ws.Cells.Item(116, 13).Value = -27470
$ws.Cells.Item(116, 14).Value = -37885.25
$ws.Cells.Item(122, 8).Value = 1863.6086
$ws.Cells.Item(122, 9).Value = 1450.7894
$ws.Cells.Item(122, 10).Value = 3824.5
$ws.Cells.Item(122, 11).Value = 4352.3682
$ws.Cells.Item(122, 12).Value = 11473.5
$ws.Cells.Item(122, 13).Value = -1902.3682
$ws.Cells.Item(122, 14).Value = -16373.5
$ws.Cells.Item(137, 8).Value = 2043.5
$ws.Cells.Item(137, 9).Value = 2065
$ws.Cells.Item(137, 11).Value = 6195
$ws.Cells.Item(137, 13).Value = -3645
$ws.Cells.Item(138, 8).Value = 10250.244
$ws.Cells.Item(138, 10).Value = 10474.646
$ws.Cells.Item(138, 12).Value = 31423.938
$ws.Cells.Item(138, 14).Value = -41703.938

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(43, 8).Value = 19599.666
$ws.Cells.Item(43, 9).Value = 11342
$ws.Cells.Item(43, 10).Value = 20631.875
$ws.Cells.Item(43, 11).Value = 11342
$ws.Cells.Item(43, 12).Value = 20631.875
$ws.Cells.Item(43, 13).Value = -11029
$ws.Cells.Item(43, 14).Value = -21257.875
$ws.Cells.Item(44, 8).Value = 31049
$ws.Cells.Item(44, 10).Value = 31049
$ws.Cells.Item(44, 12).Value = 31049
$ws.Cells.Item(44, 14).Value = -32025
$ws.Cells.Item(55, 8).Value = 25526.5
$ws.Cells.Item(55, 10).Value = 31053
$ws.Cells.Item(55, 12).Value = 31053
$ws.Cells.Item(55, 14).Value = -31683
$ws.Cells.Item(110, 8).Value = 210747.25
$ws.Cells.Item(110, 9).Value = 240158.77
$ws.Cells.Item(110, 11).Value = 240158.77
$ws.Cells.Item(110, 13).Value = -238113.77

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(35, 8).Value = 93828.57000000001
$ws.Cells.Item(35, 10).Value = 93828.57000000001
$ws.Cells.Item(35, 12).Value = 93828.57000000001
$ws.Cells.Item(35, 14).Value = -94448.57000000001
$ws.Cells.Item(114, 8).Value = 100000
$ws.Cells.Item(114, 10).Value = 100000
$ws.Cells.Item(114, 12).Value = 100000
$ws.Cells.Item(114, 14).Value = -108678
$ws.Cells.Item(134, 8).Value = 20603.46
$ws.Cells.Item(134, 9).Value = 4833.852
$ws.Cells.Item(134, 11).Value = 14501.556
$ws.Cells.Item(134, 13).Value = -11966.556

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 84561.53999999999
$ws.Cells.Item(31, 9).Value = 4481.8184
$ws.Cells.Item(31, 10).Value = 525000
$ws.Cells.Item(31, 11).Value = 4481.8184
$ws.Cells.Item(31, 12).Value = 525000
$ws.Cells.Item(31, 13).Value = -4186.8184
$ws.Cells.Item(31, 14).Value = -525590
$ws.Cells.Item(34, 8).Value = 84561.53999999999
$ws.Cells.Item(34, 9).Value = 4481.8184
$ws.Cells.Item(34, 10).Value = 525000
$ws.Cells.Item(34, 11).Value = 4481.8184
$ws.Cells.Item(34, 12).Value = 525000
$ws.Cells.Item(34, 13).Value = -4279.8184
$ws.Cells.Item(34, 14).Value = -525404

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(42, 8).Value = 400
$ws.Cells.Item(42, 9).Value = 400
$ws.Cells.Item(42, 10).Value = 0
$ws.Cells.Item(42, 11).Value = 1200
$ws.Cells.Item(42, 12).Value = 0
$ws.Cells.Item(42, 13).Value = -666
$ws.Cells.Item(42, 14).ClearContents()
$ws.Cells.Item(131, 8).Value = 172113.58
$ws.Cells.Item(131, 10).Value = 118210.445
$ws.Cells.Item(131, 12).Value = 354631.335
$ws.Cells.Item(131, 14).Value = -364711.335

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1546.5
$ws.Cells.Item(97, 9).Value = 818.8182
$ws.Cells.Item(97, 10).Value = 2435.889
$ws.Cells.Item(97, 11).Value = 818.8182
$ws.Cells.Item(97, 12).Value = 2435.889
$ws.Cells.Item(97, 13).Value = -322.8182
$ws.Cells.Item(97, 14).Value = -3427.889
$ws.Cells.Item(122, 8).Value = 4012
$ws.Cells.Item(122, 9).Value = 3332.6667
$ws.Cells.Item(122, 11).Value = 9998.000100000001
$ws.Cells.Item(122, 13).Value = -7548.000100000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 5500.7
$ws.Cells.Item(7, 9).Value = 4102
$ws.Cells.Item(7, 11).Value = 4102
$ws.Cells.Item(7, 13).Value = -3990
$ws.Cells.Item(122, 8).Value = 8909
$ws.Cells.Item(122, 9).Value = 8000
$ws.Cells.Item(122, 10).Value = 8999.9
$ws.Cells.Item(122, 11).Value = 24000
$ws.Cells.Item(122, 12).Value = 26999.7
$ws.Cells.Item(122, 13).Value = -21550
$ws.Cells.Item(122, 14).Value = -31899.7
$ws.Cells.Item(126, 8).Value = 5500.7
$ws.Cells.Item(126, 9).Value = 4102
$ws.Cells.Item(126, 11).Value = 12306
$ws.Cells.Item(126, 13).Value = -9836
$ws.Cells.Item(133, 8).Value = 50000
$ws.Cells.Item(133, 10).Value = 50000
$ws.Cells.Item(133, 12).Value = 50000
$ws.Cells.Item(133, 14).Value = -55060
$ws.Cells.Item(136, 8).Value = 11168.32
$ws.Cells.Item(136, 9).Value = 8957.786
$ws.Cells.Item(136, 10).Value = 13981.728
$ws.Cells.Item(136, 11).Value = 26873.358
$ws.Cells.Item(136, 12).Value = 41945.18399999999
$ws.Cells.Item(136, 13).Value = -24323.358
$ws.Cells.Item(136, 14).Value = -47045.18399999999

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(20, 8).Value = 6260.75
$ws.Cells.Item(20, 9).Value = 10010
$ws.Cells.Item(20, 10).Value = 5011
$ws.Cells.Item(20, 11).Value = 10010
$ws.Cells.Item(20, 12).Value = 5011
$ws.Cells.Item(20, 13).Value = -9770
$ws.Cells.Item(20, 14).Value = -5491
$ws.Cells.Item(136, 8).Value = 262457.72
$ws.Cells.Item(136, 9).Value = 339497.88
$ws.Cells.Item(136, 11).Value = 1018493.64
$ws.Cells.Item(136, 13).Value = -1015943.64
